$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "78.8/140"
